# Updates the crypto tracker sheet with the latest Price (column D) and
# 1h Volume change (column E) figures for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the refreshed prices are plain numeric-looking strings
# (e.g. "0.999"); mark those cells as Text first so Excel keeps storing
# them the same way the rest of the Price column already does, instead
# of silently converting them to numbers.
$textFormatCells = @("D4","D5","D6","D7","D9","D10","D11","D12","D19","D20","D22","D23","D25","D26","D29","D30","D32","D33","D36","D37","D38","D41","D42","D47","D50","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated prices (column D)
$ws.Range("D2").Value = "42.981.46"
$ws.Range("D3").Value = "2.576.54"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").Value = "315.37"
$ws.Range("D6").Value = "100.22"
$ws.Range("D7").Value = "0.575"
$ws.Range("D9").Value = "0.537"
$ws.Range("D10").Value = "36.30"
$ws.Range("D11").Value = "0.0814"
$ws.Range("D12").Value = "7.56"
$ws.Range("D13").Value = "2.972.59"
$ws.Range("D16").Value = "2.523.39"
$ws.Range("D18").Value = "43.015.95"
$ws.Range("D19").Value = "6.88"
$ws.Range("D20").Value = "12.71"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("D22").Value = "69.49"
$ws.Range("D23").Value = "250.33"
$ws.Range("D25").Value = "2.09"
$ws.Range("D26").Value = "27.14"
$ws.Range("D29").Value = "40.79"
$ws.Range("D30").Value = "10.33"
$ws.Range("D32").Value = "157.50"
$ws.Range("D33").Value = "3.45"
$ws.Range("D36").Value = "2.68"
$ws.Range("D37").Value = "18.76"
$ws.Range("D38").Value = "2.55"
$ws.Range("D41").Value = "23.78"
$ws.Range("D42").Value = "4.11"
$ws.Range("D46").Value = "2.005.24"
$ws.Range("D47").Value = "8.94"
$ws.Range("D48").Value = "2.823.41"
$ws.Range("D50").Value = "75.15"
$ws.Range("D51").Value = "82.02"

# Updated 1h volume change percentages (column E)
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  +2.32%  "
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  -3.98%  "
$ws.Range("E38").Value = "  +8.92%  "
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("E42").Value = "  +8.48%  "
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("E45").Value = "  -1.53%  "
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  -0.77%  "
$ws.Range("E51").Value = "  -3.07%  "
